# AutoCommit_11 апреля 2024 г. 16:52:00_SibNout2023
#
# 1) Highlight the homework-score block (C4:F30) with a solid green fill.
# 2) Make every previously-blank score cell in that block an explicit 0.
# 3) Add a new ДЗ_5 data column (L4:L30), all zero for now.
# 4) Restore the view state (active cell / selection) to what the author
#    left the sheet in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Green fill across the whole homework-grade block -------------------
$ws.Range("C4:F30").Interior.Color = 5296274   # 0x50D092 = RGB(146,208,80) = FF92D050

# --- 2) Turn blank grade cells into explicit zeros --------------------------
$blankCells = @(
    "F4",
    "C6", "D6", "E6", "F6",
    "D7", "E7", "F7",
    "C9", "D9", "E9", "F9",
    "C10", "D10", "E10", "F10",
    "E13",
    "D14", "E14", "F14",
    "C15", "D15", "E15", "F15",
    "C16", "D16", "E16", "F16",
    "E17",
    "F18",
    "C19", "D19", "E19", "F19",
    "F21",
    "D24", "E24",
    "F27",
    "C28", "D28", "E28", "F28",
    "C29", "D29", "E29", "F29",
    "D30", "F30"
)
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = 0
}

# --- 3) New ДЗ_5 data column (L), zero-filled for every student ------------
for ($r = 4; $r -le 30; $r++) {
    $ws.Cells.Item($r, 12).Value = 0
}

# --- 4) View state: scroll/selection left on L25 after the edit ------------
$ws.Range("L25").Select()
